$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.148.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.013.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.010.93"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.32"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.52%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.39"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.510.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.131.49"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.003.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "446.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.19"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.19%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.15"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.96"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0352"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.726.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.08%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.72%  "
